# Add the missing "inforequest_submitted" localization key/value pair to the
# "common" sheet, right after the existing "inforequest_answered" row.
#
# This mirrors what a translator does in Excel: select row 281 (currently
# "infoRequest_application_tab_title" / "Neuvontapyynnön tiedot"), insert a
# new blank row above it (everything below shifts down by one), and type the
# new key/value pair into the freshly inserted row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("common")
$ws.Activate()

# Insert a new row before row 281, shifting rows 281:421 down to 282:422.
$ws.Rows.Item(281).Insert()

# Fill in the newly inserted row with the new translation key and its value.
$ws.Cells.Item(281, 1).Value = "inforequest_submitted"
$ws.Cells.Item(281, 2).Value = "Vireillä"

# Match the recorded viewport/selection after the edit.
$ws.Application.ActiveWindow.ScrollRow = 247
$ws.Range("B282").Select()
